# Update leveling profit figures in the per-job Siren Profits workbook.
# Values mirror a refreshed Universalis market-price pull for the affected leves.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Cells.Item(15, 8).Value = 2223.2903
$ws.Cells.Item(15, 9).Value = 2223.2903
$ws.Cells.Item(15, 11).Value = 6669.8709
$ws.Cells.Item(15, 13).Value = -6500.8709

# Row 43: Growing Is Knowing
$ws.Cells.Item(43, 8).Value = 9910.583000000001
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 9910.583000000001
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 9910.583000000001
$ws.Cells.Item(43, 14).Value = -10048.583
$ws.Cells.Item(43, 13).ClearContents()

# Row 133: Big Brush, Big Dreams
$ws.Cells.Item(133, 8).Value = 89998
$ws.Cells.Item(133, 10).Value = 89998
$ws.Cells.Item(133, 12).Value = 89998
$ws.Cells.Item(133, 14).Value = -100118

# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 14779.6875
$ws.Cells.Item(137, 9).Value = 19952.273
$ws.Cells.Item(137, 10).Value = 3400
$ws.Cells.Item(137, 11).Value = 59856.819
$ws.Cells.Item(137, 12).Value = 10200
$ws.Cells.Item(137, 13).Value = -57306.819
$ws.Cells.Item(137, 14).Value = -15300

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks
$ws.Cells.Item(45, 8).Value = 95834.95
$ws.Cells.Item(45, 9).Value = 205460.2
$ws.Cells.Item(45, 10).Value = 4480.5835
$ws.Cells.Item(45, 11).Value = 205460.2
$ws.Cells.Item(45, 12).Value = 4480.5835
$ws.Cells.Item(45, 13).Value = -205083.2
$ws.Cells.Item(45, 14).Value = -5234.5835

# Row 61: Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 16717.5
$ws.Cells.Item(61, 9).Value = 20912.666
$ws.Cells.Item(61, 10).Value = 4132
$ws.Cells.Item(61, 11).Value = 20912.666
$ws.Cells.Item(61, 12).Value = 4132
$ws.Cells.Item(61, 13).Value = -20700.666
$ws.Cells.Item(61, 14).Value = -4556

# Row 63: Rivets Run through It
$ws.Cells.Item(63, 8).Value = 1624.75
$ws.Cells.Item(63, 10).Value = 1249
$ws.Cells.Item(63, 12).Value = 1249
$ws.Cells.Item(63, 14).Value = -2621

# Row 66: A Riveting Revival (L)
$ws.Cells.Item(66, 8).Value = 1624.75
$ws.Cells.Item(66, 10).Value = 1249
$ws.Cells.Item(66, 12).Value = 6245
$ws.Cells.Item(66, 14).Value = -13109

# Row 88: The Mast Chance
$ws.Cells.Item(88, 8).Value = 333336670
$ws.Cells.Item(88, 9).Value = 5000
$ws.Cells.Item(88, 10).Value = 500002500
$ws.Cells.Item(88, 11).Value = 5000
$ws.Cells.Item(88, 12).Value = 500002500
$ws.Cells.Item(88, 13).Value = -4594
$ws.Cells.Item(88, 14).Value = -500003312

# Row 91: The Rose and the Riveter (L)
$ws.Cells.Item(91, 8).Value = 333336670
$ws.Cells.Item(91, 9).Value = 5000
$ws.Cells.Item(91, 10).Value = 500002500
$ws.Cells.Item(91, 11).Value = 5000
$ws.Cells.Item(91, 12).Value = 500002500
$ws.Cells.Item(91, 13).Value = -3596
$ws.Cells.Item(91, 14).Value = -500005308

# Row 97: Ore for Me
$ws.Cells.Item(97, 8).Value = 9529258
$ws.Cells.Item(97, 9).Value = 8933.25
$ws.Cells.Item(97, 10).Value = 22223024
$ws.Cells.Item(97, 11).Value = 8933.25
$ws.Cells.Item(97, 12).Value = 22223024
$ws.Cells.Item(97, 13).Value = -8437.25
$ws.Cells.Item(97, 14).Value = -22224016

# Row 122: Haste for High Durium
$ws.Cells.Item(122, 8).Value = 1298673.5
$ws.Cells.Item(122, 9).Value = 4175.84
$ws.Cells.Item(122, 11).Value = 12527.52
$ws.Cells.Item(122, 13).Value = -10077.52

# Row 136: Metal with Mettle
$ws.Cells.Item(136, 8).Value = 16717.5
$ws.Cells.Item(136, 9).Value = 20912.666
$ws.Cells.Item(136, 10).Value = 4132
$ws.Cells.Item(136, 11).Value = 62737.99800000001
$ws.Cells.Item(136, 12).Value = 12396
$ws.Cells.Item(136, 13).Value = -60187.99800000001
$ws.Cells.Item(136, 14).Value = -17496

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight
$ws.Cells.Item(64, 8).Value = 7871.909
$ws.Cells.Item(64, 10).Value = 1954.5555
$ws.Cells.Item(64, 12).Value = 1954.5555
$ws.Cells.Item(64, 14).Value = -2404.5555

# Row 67: Bearing the Brunt (L)
$ws.Cells.Item(67, 8).Value = 7871.909
$ws.Cells.Item(67, 10).Value = 1954.5555
$ws.Cells.Item(67, 12).Value = 1954.5555
$ws.Cells.Item(67, 14).Value = -3514.5555

# Row 130: Annals of the Empire I
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).ClearContents()

# Row 132: Always Be Prepaired
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).ClearContents()

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 14369.583
$ws.Cells.Item(31, 9).Value = 20562.143
$ws.Cells.Item(31, 11).Value = 20562.143
$ws.Cells.Item(31, 13).Value = -20267.143

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 14369.583
$ws.Cells.Item(34, 9).Value = 20562.143
$ws.Cells.Item(34, 11).Value = 20562.143
$ws.Cells.Item(34, 13).Value = -20360.143

# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 3200.577
$ws.Cells.Item(58, 9).Value = 3115.1667
$ws.Cells.Item(58, 10).Value = 3392.75
$ws.Cells.Item(58, 11).Value = 3115.1667
$ws.Cells.Item(58, 12).Value = 3392.75
$ws.Cells.Item(58, 13).Value = -2912.1667
$ws.Cells.Item(58, 14).Value = -3798.75

# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 8866.5625
$ws.Cells.Item(122, 9).Value = 8866.5625
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 26599.6875
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -24149.6875
$ws.Cells.Item(122, 14).ClearContents()

# Row 134: Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 2389.0715
$ws.Cells.Item(134, 9).Value = 2473.739
$ws.Cells.Item(134, 10).Value = 1999.6
$ws.Cells.Item(134, 11).Value = 7421.217000000001
$ws.Cells.Item(134, 12).Value = 5998.799999999999
$ws.Cells.Item(134, 13).Value = -4886.217000000001
$ws.Cells.Item(134, 14).Value = -11068.8

# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 3200.577
$ws.Cells.Item(136, 9).Value = 3115.1667
$ws.Cells.Item(136, 10).Value = 3392.75
$ws.Cells.Item(136, 11).Value = 9345.500100000001
$ws.Cells.Item(136, 12).Value = 10178.25
$ws.Cells.Item(136, 13).Value = -6795.500100000001
$ws.Cells.Item(136, 14).Value = -15278.25

# Row 141: No Greater Treasure
$ws.Cells.Item(141, 8).Value = 343645
$ws.Cells.Item(141, 9).Value = 90000
$ws.Cells.Item(141, 10).Value = 371827.78
$ws.Cells.Item(141, 11).Value = 90000
$ws.Cells.Item(141, 12).Value = 371827.78
$ws.Cells.Item(141, 13).Value = -84820
$ws.Cells.Item(141, 14).Value = -382187.78

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 13604.833
$ws.Cells.Item(102, 9).Value = 19752.143
$ws.Cells.Item(102, 10).Value = 4998.6
$ws.Cells.Item(102, 11).Value = 19752.143
$ws.Cells.Item(102, 12).Value = 4998.6
$ws.Cells.Item(102, 13).Value = -18130.143
$ws.Cells.Item(102, 14).Value = -8242.6

# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 9021.576999999999
$ws.Cells.Item(122, 9).Value = 5773.45
$ws.Cells.Item(122, 10).Value = 19848.666
$ws.Cells.Item(122, 11).Value = 17320.35
$ws.Cells.Item(122, 12).Value = 59545.99800000001
$ws.Cells.Item(122, 13).Value = -14870.35
$ws.Cells.Item(122, 14).Value = -64445.99800000001

# Row 126: Gold Rush Order
$ws.Cells.Item(126, 8).Value = 18913.166
$ws.Cells.Item(126, 9).Value = 120000
$ws.Cells.Item(126, 10).Value = 9723.454
$ws.Cells.Item(126, 11).Value = 360000
$ws.Cells.Item(126, 12).Value = 29170.362
$ws.Cells.Item(126, 13).Value = -357530
$ws.Cells.Item(126, 14).Value = -34110.362

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 17: Only the Best
$ws.Cells.Item(17, 8).Value = 149.5
$ws.Cells.Item(17, 10).Value = 149
$ws.Cells.Item(17, 12).Value = 149
$ws.Cells.Item(17, 14).Value = -489

# Row 22: Skin off Their Backs
$ws.Cells.Item(22, 8).Value = 6568.8237
$ws.Cells.Item(22, 9).Value = 10375.6
$ws.Cells.Item(22, 10).Value = 1130.5714
$ws.Cells.Item(22, 11).Value = 10375.6
$ws.Cells.Item(22, 12).Value = 1130.5714
$ws.Cells.Item(22, 13).Value = -10080.6
$ws.Cells.Item(22, 14).Value = -1720.5714

# Row 27: Fire and Hide
$ws.Cells.Item(27, 8).Value = 6568.8237
$ws.Cells.Item(27, 9).Value = 10375.6
$ws.Cells.Item(27, 10).Value = 1130.5714
$ws.Cells.Item(27, 11).Value = 10375.6
$ws.Cells.Item(27, 12).Value = 1130.5714
$ws.Cells.Item(27, 13).Value = -10268.6
$ws.Cells.Item(27, 14).Value = -1344.5714

# Row 40: Best Served Toad
$ws.Cells.Item(40, 8).Value = 28702.666
$ws.Cells.Item(40, 9).Value = 32663.428
$ws.Cells.Item(40, 10).Value = 20781.143
$ws.Cells.Item(40, 11).Value = 32663.428
$ws.Cells.Item(40, 12).Value = 20781.143
$ws.Cells.Item(40, 13).Value = -32527.428
$ws.Cells.Item(40, 14).Value = -21053.143

# Row 42: Slave to Fashion
$ws.Cells.Item(42, 8).Value = 88400
$ws.Cells.Item(42, 9).Value = 88400
$ws.Cells.Item(42, 11).Value = 88400
$ws.Cells.Item(42, 13).Value = -87837

# Row 49: First They Came for the Heretics
$ws.Cells.Item(49, 8).Value = 88400
$ws.Cells.Item(49, 9).Value = 88400
$ws.Cells.Item(49, 11).Value = 88400
$ws.Cells.Item(49, 13).Value = -88253

# Row 122: Hell on Leather
$ws.Cells.Item(122, 8).Value = 6013.6875
$ws.Cells.Item(122, 9).Value = 8012
$ws.Cells.Item(122, 10).Value = 4459.4443
$ws.Cells.Item(122, 11).Value = 24036
$ws.Cells.Item(122, 12).Value = 13378.3329
$ws.Cells.Item(122, 13).Value = -21586
$ws.Cells.Item(122, 14).Value = -18278.3329

# Row 135: Dreams of Ja
$ws.Cells.Item(135, 8).Value = 93927.2
$ws.Cells.Item(135, 10).Value = 93927.2
$ws.Cells.Item(135, 12).Value = 93927.2
$ws.Cells.Item(135, 14).Value = -104067.2

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 5: Hire in the Blood
$ws.Cells.Item(5, 8).Value = 13500000
$ws.Cells.Item(5, 9).Value = 17333332
$ws.Cells.Item(5, 11).Value = 17333332
$ws.Cells.Item(5, 13).Value = -17333220

# Row 100: Of Great Import
$ws.Cells.Item(100, 8).Value = 36404.312
$ws.Cells.Item(100, 9).Value = 21798.75
$ws.Cells.Item(100, 11).Value = 43597.5
$ws.Cells.Item(100, 13).Value = -43056.5

# Row 122: Heavy Armoire
$ws.Cells.Item(122, 8).Value = 5082.0356
$ws.Cells.Item(122, 9).Value = 1814.9
$ws.Cells.Item(122, 10).Value = 13249.875
$ws.Cells.Item(122, 11).Value = 5444.700000000001
$ws.Cells.Item(122, 12).Value = 39749.625
$ws.Cells.Item(122, 13).Value = -2994.700000000001
$ws.Cells.Item(122, 14).Value = -44649.625
